# Insert a new weekly price record for Espinaca at Terminal La Palmera de
# La Serena. The new row is inserted above the current row 66, shifting the
# existing data rows (old 66-189) down by one (to 67-190) while preserving
# their formatting, and the freshly inserted row 66 is populated with the
# new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(66).Insert()

$ws.Cells.Item(66, 1).Value  = 8
$ws.Cells.Item(66, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(66, 3).Value  = "Coquimbo"
$ws.Cells.Item(66, 4).Value  = 44533
$ws.Cells.Item(66, 5).Value  = 4
$ws.Cells.Item(66, 6).Value  = 100112012
$ws.Cells.Item(66, 7).Value  = "Espinaca"
$ws.Cells.Item(66, 8).Value  = "Sin especificar"
$ws.Cells.Item(66, 9).Value  = "Primera"
$ws.Cells.Item(66, 10).Value = 3360
$ws.Cells.Item(66, 11).Value = 400
$ws.Cells.Item(66, 12).Value = 500
$ws.Cells.Item(66, 13).Value = 450
$ws.Cells.Item(66, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(66, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(66, 16).Value = 900
$ws.Cells.Item(66, 17).Value = 0.5
$ws.Cells.Item(66, 18).Value = "Hortaliza"
